# Update odds/score data on Sheet1 for the Flashscore weekly games workbook.
# The workbook's ActiveWorkbook / ActiveSheet are already open in $excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 3.75
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 1.91
$ws.Range("J2").Value = 4.33
$ws.Range("L2").Value = 2.6
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.5
$ws.Range("Z2").Value = 19
$ws.Range("AC2").Value = 29
$ws.Range("AD2").Value = 34
$ws.Range("AF2").Value = 6.5
$ws.Range("AJ2").Value = 9.5
$ws.Range("AL2").Value = 17

# --- Row 3 updates ---
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 3.8
$ws.Range("K3").Value = 2.3
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.83
$ws.Range("R3").Value = 2.03
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 1.36
$ws.Range("U3").Value = 1.36
$ws.Range("V3").Value = 3
$ws.Range("W3").Value = 1.83
$ws.Range("X3").Value = 1.83
$ws.Range("Y3").Value = 15
$ws.Range("AB3").Value = 67
$ws.Range("AE3").Value = 11
$ws.Range("AG3").Value = 17
$ws.Range("AI3").Value = 7
$ws.Range("AJ3").Value = 7.5
$ws.Range("AM3").Value = 13
$ws.Range("AN3").Value = 26
$ws.Range("AO3").Value = 251

# --- Row 6 updates ---
$ws.Range("Q6").Value = 1.75
$ws.Range("R6").Value = 2.05

# --- Row 9 updates ---
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 1.73
